# Weekly update: insert a new price record for "Berenjena" (Vega Modelo de
# Temuco) as row 438, shifting the existing rows 438-467 down to 439-468.
# The sheet's used range grows from A1:R467 to A1:R468.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 438..467 down one position (Excel recalculates the dimension
# / used range automatically on insert).
$ws.Rows.Item(438).Insert()

# Populate the newly inserted row 438 with this week's entry. Columns
# A, B, C, E, F, G, H, I, O, P, R keep the same values as the row that used
# to occupy position 438 (now 439); only the date/volume/price/unit/qty
# columns (D, J, K, L, M, N, Q) carry new data.
$ws.Range("A438").Value = 10
$ws.Range("B438").Value = "Vega Modelo de Temuco"
$ws.Range("C438").Value = "La Araucanía"
$ws.Range("D438").Value = 45106
$ws.Range("E438").Value = 9
$ws.Range("F438").Value = 100112001
$ws.Range("G438").Value = "Berenjena"
$ws.Range("H438").Value = "Sin especificar"
$ws.Range("I438").Value = "Primera"
$ws.Range("J438").Value = 35
$ws.Range("K438").Value = 10000
$ws.Range("L438").Value = 10000
$ws.Range("M438").Value = 10000
$ws.Range("N438").Value = "$/caja 40 unidades"
$ws.Range("O438").Value = "Región de Arica y Parinacota"
$ws.Range("P438").Value = 250
$ws.Range("Q438").Value = 40
$ws.Range("R438").Value = "Hortaliza"
